$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The loading-time export script now appends one more day's worth of
# columns (captured while the automation run was catching/reporting
# exceptions) instead of stopping at the last successful date column.
# New date label " Oct 19" goes across the header row, and the new
# per-category counts land in the matching data rows below it.

# Header row: new date column repeated across G1:V1
$ws.Range("G1:V1").Style = "Normal"
$ws.Range("G1:V1").Value = " Oct 19"

# Row 2 (Login)
$ws.Range("G2:P2").Style = "Normal"
$ws.Range("G2:P2").Value = 0.0
$ws.Range("T2:V2").Style = "Normal"
$ws.Range("T2:V2").Value = 0.0

# Row 3 (Dashboard)
$ws.Range("G3:P3").Style = "Normal"
$ws.Range("G3").Value = 8.0
$ws.Range("H3").Value = 13.0
$ws.Range("I3").Value = 9.0
$ws.Range("J3").Value = 8.0
$ws.Range("K3").Value = 9.0
$ws.Range("L3").Value = 8.0
$ws.Range("M3").Value = 8.0
$ws.Range("N3").Value = 9.0
$ws.Range("O3").Value = 9.0
$ws.Range("P3").Value = 9.0
$ws.Range("T3:V3").Style = "Normal"
$ws.Range("T3").Value = 7.0
$ws.Range("U3").Value = 9.0
$ws.Range("V3").Value = 9.0

# Row 4 (Trucks)
$ws.Range("G4:V4").Style = "Normal"
$ws.Range("G4").Value = 27.0
$ws.Range("H4").Value = 35.0
$ws.Range("I4").Value = 25.0
$ws.Range("J4").Value = 28.0
$ws.Range("K4").Value = 30.0
$ws.Range("L4").Value = 26.0
$ws.Range("M4").Value = 32.0
$ws.Range("N4").Value = 25.0
$ws.Range("O4").Value = 39.0
$ws.Range("P4").Value = 27.0
$ws.Range("Q4").Value = 24.0
$ws.Range("R4").Value = 24.0
$ws.Range("S4").Value = 34.0
$ws.Range("T4").Value = 46.0
$ws.Range("U4").Value = 38.0
$ws.Range("V4").Value = 28.0

# Row 5 (Deleted Trailers)
$ws.Range("K5:M5").Style = "Normal"
$ws.Range("K5:M5").Value = 0.0
$ws.Range("O5:P5").Style = "Normal"
$ws.Range("O5:P5").Value = 0.0
$ws.Range("U5:V5").Style = "Normal"
$ws.Range("U5:V5").Value = 0.0

# Row 6 (Deleted Trucks)
$ws.Range("K6:M6").Style = "Normal"
$ws.Range("K6:M6").Value = 0.0
$ws.Range("O6:P6").Style = "Normal"
$ws.Range("O6:P6").Value = 0.0
$ws.Range("T6:V6").Style = "Normal"
$ws.Range("T6:V6").Value = 0.0

# Row 7 (Trailer)
$ws.Range("G7:P7").Style = "Normal"
$ws.Range("G7:P7").Value = 0.0
$ws.Range("U7:V7").Style = "Normal"
$ws.Range("U7:V7").Value = 0.0
